# Updates the "Estado de Cuenta" worker table (rows 16-37): removes the
# previous account-statement rows and replaces them with the new set of
# workers / periods, per the commit "Elimna EC anteriores y se agregan
# nuevos, se modifica base de datos".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# r -> (DocNumber, Name, Period, ValorMora)
$data = @{
    16 = @("1100334373", "EDUARDO SANEZ ANGEL", "1711", 737717)
    17 = @("1100334373", "EDUARDO SANEZ ANGEL", "1710", 737717)
    18 = @("73190652", "JOSE LUIS BLANCO CABARCAS", "1711", 737717)
    19 = @("73190652", "JOSE LUIS BLANCO CABARCAS", "1710", 737717)
    20 = @("1047432610", "EDWIN RAFAEL MUÑOZ DE ANGEL", "1711", 737717)
    21 = @("1047432610", "EDWIN RAFAEL MUÑOZ DE ANGEL", "1710", 737717)
    22 = @("1052981975", "FRANKY ALI RIOS MUÑOZ", "1711", 737717)
    23 = @("1052981975", "FRANKY ALI RIOS MUÑOZ", "1710", 737717)
    24 = @("1052992439", "MIGUEL ANGEL JIMENEZ MARQUEZ", "1711", 737717)
    25 = @("1052992439", "MIGUEL ANGEL JIMENEZ MARQUEZ", "1710", 737717)
    26 = @("1052968659", "DIONISIO JOSE RIOS MUÑOZ", "1711", 737717)
    27 = @("1052968659", "DIONISIO JOSE RIOS MUÑOZ", "1710", 737717)
    28 = @("1052960660", "RAFAEL OVIDIO RIOS MUÑOZ", "1711", 737717)
    29 = @("1052960660", "RAFAEL OVIDIO RIOS MUÑOZ", "1710", 737717)
    30 = @("1052946538", "EVERTO MANUEL DE LAS OSSA ATENCIA", "1711", 737717)
    31 = @("1052946538", "EVERTO MANUEL DE LAS OSSA ATENCIA", "1710", 737717)
    32 = @("1101445961", "PEDRO LUIS TORRES ZUÑIGA", "1711", 737717)
    33 = @("1101445961", "PEDRO LUIS TORRES ZUÑIGA", "1710", 737717)
    34 = @("1052956555", "JUAN FRANCISCO MEZA LOPEZ", "1711", 737717)
    35 = @("1052956555", "JUAN FRANCISCO MEZA LOPEZ", "1710", 737717)
    36 = @("1103220109", "DIEGO ARMANDO WILCHEZ NOVOA", "1710", 737717)
    37 = @("1103220109", "DIEGO ARMANDO WILCHEZ NOVOA", "1710", 737717)
}

foreach ($r in 16..37) {
    $row = $data[$r]
    $ws.Cells.Item($r, 3).Value2 = $row[0]   # C: N Doc Trabajador
    $ws.Cells.Item($r, 4).Value2 = $row[1]   # D: Nombre Trabajador
    $ws.Cells.Item($r, 5).Value2 = $row[2]   # E: Periodo Mora
    $ws.Cells.Item($r, 7).Value2 = $row[3]   # G: Valor Mora
}

# Column widths follow the new (longer) best-fit content.
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null
$ws.Columns.Item(4).AutoFit() | Out-Null
$ws.Columns.Item(5).AutoFit() | Out-Null
$ws.Columns.Item(6).AutoFit() | Out-Null
$ws.Columns.Item(7).AutoFit() | Out-Null
$ws.Columns.Item(8).AutoFit() | Out-Null
$ws.Columns.Item(9).AutoFit() | Out-Null
$ws.Columns.Item(10).AutoFit() | Out-Null
